# Update the "Förändrad" (Changed) date column C for all data rows (2-120)
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 120; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value = 45224
    }
}
